$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D, E, G contain numeric-looking text (prices, percentages, hour)
# that Excel would auto-convert to numbers on assignment. Force each such
# target cell to Text format first so the literal string is preserved,
# exactly mirroring the original inline-string cell contents.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "294.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.82%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "12"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.69%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "12"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.014"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.62%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "12"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07393"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-4.11%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "12"

# Row 6
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.571"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-4.76%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "12"

# Row 7
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9194"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.47%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "12"

# Row 8
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.1190"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.41%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "12"

# Row 9
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1754"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.71%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "12"

# Row 10
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08697"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.93%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "12"

# Row 11
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04161"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.72%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "12"

# Row 12
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1054"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.33%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "12"

# Row 13
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001273"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.48%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "12"

# Row 14
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.005803"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.88%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "12"

# Row 15
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.414"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.87%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "12"

# Row 16
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.300"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.36%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "12"

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.11%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "12"

# Row 18
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "12"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.563"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.29%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "12"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1354"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.41%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "12"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2803"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.60%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "12"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.03842"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-4.35%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "12"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001282"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.15%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "12"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003904"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-4.56%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "12"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001292"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.79%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "12"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003727"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-95.04%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "12"

# Row 27
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "12"

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "12"

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "12"

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "12"

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "12"

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "12"

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "12"

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "12"

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "12"

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "12"

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "12"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02310"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-8.80%"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "12"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05037"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.45%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "12"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007713"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.67%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "12"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.004506"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "143.19%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "12"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1274"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.95%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "12"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007396"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "11.39%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "12"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006956"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.88%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "12"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3205"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.40%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "12"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006475"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.54%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "12"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.02%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "12"

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "11.63%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "12"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004205"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "35.50%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "12"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.02%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "12"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.02%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "12"
